$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
